$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.457.27'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.844.22'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.09'
$ws.Range("E5").Value = '  -3.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5214'
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3271'
$ws.Range("E8").Value = '  -3.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06802'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -5.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7808'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07758'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '1.849.23'
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.08'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.019'
$ws.Range("E15").Value = '  -2.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9989'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.94'
$ws.Range("E17").Value = '  -3.26%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007979'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '26.496.70'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").Value = '2.072.56'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.624'
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.573'
$ws.Range("E23").Value = '  -3.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.986'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.68'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.175'
$ws.Range("E26").Value = '  -7.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.645'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.01'
$ws.Range("E28").Value = '  -0.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.99'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.158'
$ws.Range("E30").Value = '  -3.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.133'
$ws.Range("E31").Value = '  -3.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08706'
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04839'
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7247'
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.132'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.838'
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.108'
$ws.Range("E37").Value = '  -3.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.234'
$ws.Range("E38").Value = '  -4.13%  '
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4865'
$ws.Range("E40").Value = '  -4.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9143'
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '111.32'
$ws.Range("E42").Value = '  -4.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.076'
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.756'
$ws.Range("E45").Value = '  -3.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4185'
$ws.Range("E46").Value = '  -5.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05931'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.063'
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.15'
$ws.Range("E49").Value = '  -2.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1237'
$ws.Range("E50").Value = '  -6.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8891'
$ws.Range("E51").Value = '  +1.19%  '
